$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data as captured in the commit diff.
# Force text number-format first so numeric-looking strings (e.g. "80.08")
# are not silently reinterpreted as numbers by Excel, then restore the
# default "Normal" style so no stray style index is left on the cell.
$cellUpdates = @{
    "D2" = "69.407.72"
    "E2" = "  -2.39%  "
    "D3" = "3.697.77"
    "E3" = "  -2.98%  "
    "E4" = "  +0.10%  "
    "D5" = "692.69"
    "E5" = "  -0.79%  "
    "D6" = "162.91"
    "E6" = "  -5.20%  "
    "D7" = "3.695.69"
    "E7" = "  -3.02%  "
    "E8" = "  +0.07%  "
    "E9" = "  -4.62%  "
    "E10" = "  -8.15%  "
    "D11" = "7.41"
    "E11" = "  -1.52%  "
    "D12" = "0.443"
    "E12" = "  -4.45%  "
    "D13" = "0.0000240"
    "E13" = "  -5.31%  "
    "E14" = "  -7.13%  "
    "D15" = "4.320.12"
    "E15" = "  -2.97%  "
    "D16" = "3.701.78"
    "E16" = "  -4.15%  "
    "D17" = "69.437.32"
    "E17" = "  -2.31%  "
    "E18" = "  -0.92%  "
    "D19" = "16.25"
    "E19" = "  -6.97%  "
    "E20" = "  -7.60%  "
    "D21" = "480.89"
    "E21" = "  -6.18%  "
    "D22" = "9.97"
    "E22" = "  -6.63%  "
    "E23" = "  -7.22%  "
    "D24" = "80.08"
    "D25" = "3.845.24"
    "E25" = "  -2.98%  "
    "E26" = "  -9.17%  "
    "E27" = "  +0.02%  "
    "D28" = "11.39"
    "E28" = "  -5.57%  "
    "E29" = "  -8.54%  "
    "E30" = "  -10.53%  "
    "D31" = "2.73"
    "E31" = "  -9.80%  "
    "D32" = "6.84"
    "E32" = "  -7.89%  "
    "E33" = "  -7.82%  "
    "E34" = "  -4.49%  "
    "D35" = "27.02"
    "E35" = "  -7.13%  "
    "E36" = "  +0.27%  "
    "D37" = "3.665.47"
    "E37" = "  -2.81%  "
    "D38" = "8.51"
    "E38" = "  -7.13%  "
    "D39" = "6.37"
    "E39" = "  +6.18%  "
    "E40" = "  -2.84%  "
    "E43" = "  +0.01%  "
    "D44" = "0.955"
    "E44" = "  -6.41%  "
    "D45" = "163.89"
    "E45" = "  -5.26%  "
    "D46" = "47.91"
    "E46" = "  -3.06%  "
    "D47" = "30.26"
    "E47" = "  +2.49%  "
    "E48" = "  -15.04%  "
    "E49" = "  -1.34%  "
    "D50" = "1.35"
    "E50" = "  -1.37%  "
    "D51" = "0.000286"
    "E51" = "  -7.81%  "
}

foreach ($addr in $cellUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $cellUpdates[$addr]
    $cell.Style = "Normal"
}
